# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# OFF sheet - update Road (row 3) target depth totals
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 403
$wsOff.Range("C3").Value = 287
$wsOff.Range("D3").Value = 98
$wsOff.Range("E3").Value = 48
$wsOff.Range("F3").Value = 7

# DEF sheet - update Road (row 3) target depth totals
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 427
$wsDef.Range("C3").Value = 298
$wsDef.Range("D3").Value = 102
$wsDef.Range("E3").Value = 54
